$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: build a paragraph out of N separate text "runs" without letting
# same-formatted adjacent runs get auto-merged. We do this by inserting each
# run's text into its OWN paragraph (so they can't merge), then deleting the
# paragraph marks between them to splice the paragraphs back into one.
# ---------------------------------------------------------------------------

# ===========================================================================
# 1) Last existing paragraph: split "...and compare to the Cessna..." into
#    three runs, changing "compare" -> "compare it".
# ===========================================================================
$lastPara = $d.Paragraphs.Last
$rng = $lastPara.Range
$rng.Text = ""
$rng.InsertAfter("I’ll likely use VLM over a range of angles of attack and ")

$tmp = $rng.InsertParagraphAfter()
$p2 = $d.Paragraphs.Last
$p2.Range.InsertAfter("compare it")

$tmp = $p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Last
$p3.Range.InsertAfter(" to the Cessna 172 lift curves to see how accurate my method is. After that I will optimize.")

# Splice the three paragraphs back together by deleting the paragraph marks
# that separate them (this preserves each chunk as its own <w:r>).
$mark1 = $d.Range($lastPara.Range.End - 1, $lastPara.Range.End)
$mark1.Delete()
$mark2 = $d.Range($lastPara.Range.End - 1, $lastPara.Range.End)
$mark2.Delete()

# ===========================================================================
# 2) New paragraph: "Use web plot digitizer to interpolate any graphs that
#    I'm pulling from research."
# ===========================================================================
$rng = $lastPara.Range
$rng.Collapse(0)
$rng.InsertParagraphAfter()
$pWeb = $d.Paragraphs.Last
$pWeb.Range.InsertAfter("Use web plot digitizer to interpolate any graphs that I’m pulling from research.")

# ===========================================================================
# 3) New empty paragraph.
# ===========================================================================
$rng = $pWeb.Range
$rng.Collapse(0)
$rng.InsertParagraphAfter()
$pEmpty = $d.Paragraphs.Last

# ===========================================================================
# 4) New bold paragraph with three runs (all bold):
#    'Comparing results to the study done in the paper "' +
#    'Comparative Study and Aerodynamic Analysis of Rectangular Wing Using
#     High-Lift Systems' + '."'
# ===========================================================================
$rng = $pEmpty.Range
$rng.Collapse(0)
$rng.InsertParagraphAfter()
$bp1 = $d.Paragraphs.Last
$bp1.Range.Font.Bold = 1
$bp1.Range.Font.BoldBi = 1
$bp1.Range.InsertAfter("Comparing results to the study done in the paper “")

$tmp = $bp1.Range.InsertParagraphAfter()
$bp2 = $d.Paragraphs.Last
$bp2.Range.Font.Bold = 1
$bp2.Range.Font.BoldBi = 1
$bp2.Range.InsertAfter("Comparative Study and Aerodynamic Analysis of Rectangular Wing Using High-Lift Systems")

$tmp = $bp2.Range.InsertParagraphAfter()
$bp3 = $d.Paragraphs.Last
$bp3.Range.Font.Bold = 1
$bp3.Range.Font.BoldBi = 1
$bp3.Range.InsertAfter(".”")

$markb1 = $d.Range($bp1.Range.End - 1, $bp1.Range.End)
$markb1.Delete()
$markb2 = $d.Range($bp1.Range.End - 1, $bp1.Range.End)
$markb2.Delete()

# ===========================================================================
# 5) New paragraph: "For figure 23 the lift coefficient at 2 degrees angle of
#    attack is 0.519." (NOT bold)
# ===========================================================================
$rng = $bp1.Range
$rng.Collapse(0)
$rng.InsertParagraphAfter()
$pFig = $d.Paragraphs.Last
$pFig.Range.Font.Bold = 0
$pFig.Range.Font.BoldBi = 0
$pFig.Range.InsertAfter("For figure 23 the lift coefficient at 2 degrees angle of attack is 0.519.")

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
